# Weekly update: add a new date's worth of price records (4 rows) for
# "Hortaliza, Mercado Mayorista Lo Valledor de Santiago - Zanahoria".
#
# The new rows are inserted right after the existing row 731, which pushes
# all of the old rows 732:808 down to 736:812 (and grows the sheet's used
# range from A1:R808 to A1:R812).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before the current row 732, shifting the remaining
# data (previously rows 732:808) down to rows 736:812.
$ws.Rows("732:735").Insert()

# Values for the newly inserted rows (date serial 44449 = 2021-09-10).
$newRows = @(
    @{ Row = 732; Fecha = 44449; Calidad = "Primera"; Volumen = 2200; PrecioMin = 5000; PrecioMax = 5500; PrecioProm = 5182; Origen = "Chillán";              PrecioKg = 259 },
    @{ Row = 733; Fecha = 44449; Calidad = "Primera"; Volumen = 1640; PrecioMin = 4500; PrecioMax = 5000; PrecioProm = 4771; Origen = "Región Metropolitana"; PrecioKg = 239 },
    @{ Row = 734; Fecha = 44449; Calidad = "Segunda"; Volumen = 360;  PrecioMin = 4000; PrecioMax = 4000; PrecioProm = 4000; Origen = "Chillán";              PrecioKg = 200 },
    @{ Row = 735; Fecha = 44449; Calidad = "Segunda"; Volumen = 440;  PrecioMin = 3500; PrecioMax = 4000; PrecioProm = 3773; Origen = "Región Metropolitana"; PrecioKg = 189 }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $ws.Range("A$r").Value = 6
    $ws.Range("B$r").Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Range("C$r").Value = "Metropolitana"
    $ws.Range("D$r").Value = $item.Fecha
    $ws.Range("E$r").Value = 13
    $ws.Range("F$r").Value = 100114013
    $ws.Range("G$r").Value = "Zanahoria"
    $ws.Range("H$r").Value = "Sin especificar"
    $ws.Range("I$r").Value = $item.Calidad
    $ws.Range("J$r").Value = $item.Volumen
    $ws.Range("K$r").Value = $item.PrecioMin
    $ws.Range("L$r").Value = $item.PrecioMax
    $ws.Range("M$r").Value = $item.PrecioProm
    $ws.Range("N$r").Value = "`$/saco 20 kilos"
    $ws.Range("O$r").Value = $item.Origen
    $ws.Range("P$r").Value = $item.PrecioKg
    $ws.Range("Q$r").Value = 20
    $ws.Range("R$r").Value = "Hortaliza"
}
